# Update "want to go" counts (column F) in several rows across sheets.
$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 275
$ws1.Range("F3").Value = 591
$ws1.Range("F4").Value = 6926
$ws1.Range("F9").Value = 1145
$ws1.Range("F10").Value = 16509
$ws1.Range("F11").Value = 15
$ws1.Range("F17").Value = 11518
$ws1.Range("F18").Value = 22
$ws1.Range("F19").Value = 1170
$ws1.Range("F20").Value = 4551
$ws1.Range("F21").Value = 389
$ws1.Range("F22").Value = 397
$ws1.Range("F24").Value = 866
$ws1.Range("F25").Value = 329

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 3

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 275
$ws4.Range("F3").Value = 591
$ws4.Range("F4").Value = 6926
$ws4.Range("F10").Value = 1145
$ws4.Range("F11").Value = 16509
$ws4.Range("F12").Value = 15
$ws4.Range("F18").Value = 3
$ws4.Range("F20").Value = 11518
$ws4.Range("F21").Value = 22
$ws4.Range("F22").Value = 1170
$ws4.Range("F23").Value = 4551
$ws4.Range("F24").Value = 389
$ws4.Range("F25").Value = 397
$ws4.Range("F27").Value = 866
$ws4.Range("F28").Value = 329
